$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.118015681109569
$ws.Range("C2").Value = 0.3082411197234762
$ws.Range("D2").Value = 0.03355238475236177
$ws.Range("E2").Value = 0.1233058829017886
$ws.Range("F2").Value = 0.8231263785663572
$ws.Range("L2").Value = 0.1891530951565983
$ws.Range("M2").Value = 0.2361424816765876
$ws.Range("O2").Value = 2.871797905458976
# Row 3
$ws.Range("B3").Value = 1.012251070564901
$ws.Range("C3").Value = 0.294607476841918
$ws.Range("D3").Value = 0.0320621126323104
$ws.Range("E3").Value = 0.1247236308366461
$ws.Range("F3").Value = 0.8238021451527047
$ws.Range("L3").Value = 0.1865257311867694
$ws.Range("M3").Value = 0.2199033853709338
$ws.Range("O3").Value = 2.890035906224625
# Row 4
$ws.Range("B4").Value = 0.9473682001213888
$ws.Range("C4").Value = 0.2862065952472221
$ws.Range("D4").Value = 0.03114326073572471
$ws.Range("E4").Value = 0.1256438214809056
$ws.Range("F4").Value = 0.8248677082826816
$ws.Range("L4").Value = 0.1850045961093159
$ws.Range("M4").Value = 0.2099848639771267
$ws.Range("O4").Value = 2.903466340057037
# Row 5
$ws.Range("B5").Value = 0.9209436966749536
$ws.Range("C5").Value = 0.2827759353475585
$ws.Range("D5").Value = 0.03076788752337478
$ws.Range("E5").Value = 0.1260313135664708
$ws.Range("F5").Value = 0.8254652868235297
$ws.Range("L5").Value = 0.1844079178928695
$ws.Range("M5").Value = 0.2059563668816793
$ws.Range("O5").Value = 2.90949973733305
# Row 6
$ws.Range("B6").Value = 0.9165569249605028
$ws.Range("C6").Value = 0.2822058475796325
$ws.Range("D6").Value = 0.03070550148983386
$ws.Range("E6").Value = 0.1260964120892646
$ws.Range("F6").Value = 0.8255743728934704
$ws.Range("L6").Value = 0.1843102422323639
$ws.Range("M6").Value = 0.2052882526239301
$ws.Range("O6").Value = 2.910535397163414
# Row 7
$ws.Range("B7").Value = 0.9470117641444062
$ws.Range("C7").Value = 0.2861603571133458
$ws.Range("D7").Value = 0.03113820206625206
$ws.Range("E7").Value = 0.1256489966779102
$ws.Range("F7").Value = 0.824875106367891
$ws.Range("L7").Value = 0.1849964551177052
$ws.Range("M7").Value = 0.2099304797878929
$ws.Range("O7").Value = 2.903545440816259
# Row 8
$ws.Range("B8").Value = 1.081537041416539
$ws.Range("C8").Value = 0.30354656535917
$ws.Range("D8").Value = 0.03303934697024857
$ws.Range("E8").Value = 0.1237844173813557
$ws.Range("F8").Value = 0.823224186531661
$ws.Range("L8").Value = 0.1882280939914267
$ws.Range("M8").Value = 0.2305325005614094
$ws.Range("O8").Value = 2.877622509305382
# Row 9
$ws.Range("B9").Value = 1.345743109301736
$ws.Range("C9").Value = 0.3373952499985364
$ws.Range("D9").Value = 0.03673620094447472
$ws.Range("E9").Value = 0.1205216274595011
$ws.Range("F9").Value = 0.8251620350377991
$ws.Range("L9").Value = 0.1952945804908595
$ws.Range("M9").Value = 0.2713410357480726
$ws.Range("O9").Value = 2.84454345519012
# Row 10
$ws.Range("B10").Value = 1.54005357782296
$ws.Range("C10").Value = 0.3621040387107257
$ws.Range("D10").Value = 0.03943215737568551
$ws.Range("E10").Value = 0.118363539648736
$ws.Range("F10").Value = 0.8297602121856329
$ws.Range("L10").Value = 0.2009299290898809
$ws.Range("M10").Value = 0.3015652178322767
$ws.Range("O10").Value = 2.831128854499298
# Row 11
$ws.Range("B11").Value = 1.628484986013461
$ws.Range("C11").Value = 0.3733080582368871
$ws.Range("D11").Value = 0.04065404735842293
$ws.Range("E11").Value = 0.1174334934327754
$ws.Range("F11").Value = 0.8325456133071754
$ws.Range("L11").Value = 0.2035897711505328
$ws.Range("M11").Value = 0.3153663956591828
$ws.Range("O11").Value = 2.827404545917886
# Row 12
$ws.Range("B12").Value = 1.661976073157689
$ws.Range("C12").Value = 0.3775453212945195
$ws.Range("D12").Value = 0.04111607394082029
$ws.Range("E12").Value = 0.1170887276307433
$ws.Range("F12").Value = 0.8337004295918007
$ws.Range("L12").Value = 0.2046108018383705
$ws.Range("M12").Value = 0.3205998604777349
$ws.Range("O12").Value = 2.826337285004115
# Row 13
$ws.Range("B13").Value = 1.654763012999069
$ws.Range("C13").Value = 0.3766329976445206
$ws.Range("D13").Value = 0.04101659877734676
$ws.Range("E13").Value = 0.1171626491358688
$ws.Range("F13").Value = 0.8334472649933815
$ws.Range("L13").Value = 0.2043902912208893
$ws.Range("M13").Value = 0.319472420948351
$ws.Range("O13").Value = 2.826551863985912
# Row 14
$ws.Range("B14").Value = 1.631240245745119
$ws.Range("C14").Value = 0.3736567712431622
$ws.Range("D14").Value = 0.04069207228076976
$ws.Range("E14").Value = 0.1174049806346287
$ws.Range("F14").Value = 0.8326386139872142
$ws.Range("L14").Value = 0.2036734954946837
$ws.Range("M14").Value = 0.3157968117782346
$ws.Range("O14").Value = 2.827309860172477
# Row 15
$ws.Range("B15").Value = 1.616832351879282
$ws.Range("C15").Value = 0.3718330284405624
$ws.Range("D15").Value = 0.04049320150236468
$ws.Range("E15").Value = 0.1175543820454366
$ws.Range("F15").Value = 0.8321563295279617
$ws.Range("L15").Value = 0.2032362341904417
$ws.Range("M15").Value = 0.3135463336099846
$ws.Range("O15").Value = 2.827818862790878
# Row 16
$ws.Range("B16").Value = 1.534275023904343
$ws.Range("C16").Value = 0.3613710802920878
$ws.Range("D16").Value = 0.03935221079691331
$ws.Range("E16").Value = 0.1184253595077372
$ws.Range("F16").Value = 0.8295921634483108
$ws.Range("L16").Value = 0.2007580363465706
$ws.Range("M16").Value = 0.3006643063454817
$ws.Range("O16").Value = 2.831420192909547
# Row 17
$ws.Range("B17").Value = 1.483637578444302
$ws.Range("C17").Value = 0.3549435734134079
$ws.Range("D17").Value = 0.0386510740994126
$ws.Range("E17").Value = 0.1189729058225016
$ws.Range("F17").Value = 0.8281970072590354
$ws.Range("L17").Value = 0.1992623763262884
$ws.Range("M17").Value = 0.2927747686860229
$ws.Range("O17").Value = 2.834239348906266
# Row 18
$ws.Range("B18").Value = 1.454515982280725
$ws.Range("C18").Value = 0.3512432558491128
$ws.Range("D18").Value = 0.03824737554740665
$ws.Range("E18").Value = 0.1192927047498882
$ws.Range("F18").Value = 0.8274598239411546
$ws.Range("L18").Value = 0.1984111782784197
$ws.Range("M18").Value = 0.288241827455046
$ws.Range("O18").Value = 2.836084622194505
# Row 19
$ws.Range("B19").Value = 1.444656596184188
$ws.Range("C19").Value = 0.349989818402662
$ws.Range("D19").Value = 0.03811061830132445
$ws.Range("E19").Value = 0.1194018189714446
$ws.Range("F19").Value = 0.8272214280706862
$ws.Range("L19").Value = 0.1981245355897698
$ws.Range("M19").Value = 0.2867079020250998
$ws.Range("O19").Value = 2.836747799292993
# Row 20
$ws.Range("B20").Value = 1.489027646677755
$ws.Range("C20").Value = 0.3556281447870617
$ws.Range("D20").Value = 0.03872575530484568
$ws.Range("E20").Value = 0.1189141151614541
$ws.Range("F20").Value = 0.8283387662841335
$ws.Range("L20").Value = 0.1994206538836778
$ws.Range("M20").Value = 0.2936141166912734
$ws.Range("O20").Value = 2.833916077721398
# Row 21
$ws.Range("B21").Value = 1.638149355286259
$ws.Range("C21").Value = 0.374531111213031
$ws.Range("D21").Value = 0.04078741216054738
$ws.Range("E21").Value = 0.1173336006240877
$ws.Range("F21").Value = 0.8328734168321859
$ws.Range("L21").Value = 0.2038836614006243
$ws.Range("M21").Value = 0.3168762320905287
$ws.Range("O21").Value = 2.827077899255357
# Row 22
$ws.Range("B22").Value = 1.735631820369235
$ws.Range("C22").Value = 0.3868533724961196
$ws.Range("D22").Value = 0.04213086761837559
$ws.Range("E22").Value = 0.1163439030910345
$ws.Range("F22").Value = 0.8364203502453336
$ws.Range("L22").Value = 0.2068809493718931
$ws.Range("M22").Value = 0.3321215728849651
$ws.Range("O22").Value = 2.824608838518515
# Row 23
$ws.Range("B23").Value = 1.683602005027183
$ws.Range("C23").Value = 0.3802797555389645
$ws.Range("D23").Value = 0.04141421117680011
$ws.Range("E23").Value = 0.1168681678770399
$ws.Range("F23").Value = 0.8344738176604096
$ws.Range("L23").Value = 0.2052738921474742
$ws.Range("M23").Value = 0.3239810600089186
$ws.Range("O23").Value = 2.82574325534469
# Row 24
$ws.Range("B24").Value = 1.486590826166776
$ws.Range("C24").Value = 0.3553186658066068
$ws.Range("D24").Value = 0.03869199381626487
$ws.Range("E24").Value = 0.118940678821307
$ws.Range("F24").Value = 0.8282744748572952
$ws.Range("L24").Value = 0.1993490695702462
$ws.Range("M24").Value = 0.2932346385362763
$ws.Range("O24").Value = 2.834061529383831
# Row 25
$ws.Range("B25").Value = 1.274230211854615
$ws.Range("C25").Value = 0.328265636503005
$ws.Range("D25").Value = 0.03573956315217686
$ws.Range("E25").Value = 0.1213622248907576
$ws.Range("F25").Value = 0.8240817282579798
$ws.Range("L25").Value = 0.1933049387019565
$ws.Range("M25").Value = 0.2602582533825242
$ws.Range("O25").Value = 2.851584689225348
